# Apply the commit's data cleanup to the "result" sheet:
#  - Rows 192-208 get updated A/B/C/D values ("cleaning constants").
#  - Row 209 (the old last data row) is removed entirely, which also
#    shrinks the used range from A1:D209 down to A1:D208.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(192, 1).Value = 39.06
$ws.Cells.Item(192, 2).Value = 0.04000000000000001
$ws.Cells.Item(192, 3).Value = 3.896192789077759
$ws.Cells.Item(192, 4).Value = 1.09375

$ws.Cells.Item(193, 1).Value = 39.78
$ws.Cells.Item(193, 2).Value = 0.16
$ws.Cells.Item(193, 3).Value = 10.04770517349243
$ws.Cells.Item(193, 4).Value = 7.53125

$ws.Cells.Item(194, 1).Value = 39.78
$ws.Cells.Item(194, 2).Value = 0.16
$ws.Cells.Item(194, 3).Value = 23.37266087532043
$ws.Cells.Item(194, 4).Value = 14.328125

$ws.Cells.Item(195, 1).Value = 39.78
$ws.Cells.Item(195, 2).Value = 0.16
$ws.Cells.Item(195, 3).Value = 23.50989937782288
$ws.Cells.Item(195, 4).Value = 14.453125

$ws.Cells.Item(196, 1).Value = 54.5
$ws.Cells.Item(196, 2).Value = 0.16
$ws.Cells.Item(196, 3).Value = 30.71380376815796
$ws.Cells.Item(196, 4).Value = 18.90625

$ws.Cells.Item(197, 1).Value = 54.5
$ws.Cells.Item(197, 2).Value = 0.16
$ws.Cells.Item(197, 3).Value = 39.19084405899048
$ws.Cells.Item(197, 4).Value = 17.96875

$ws.Cells.Item(198, 1).Value = 39.06
$ws.Cells.Item(198, 2).Value = 0.04000000000000001
$ws.Cells.Item(198, 3).Value = 6.48872709274292
$ws.Cells.Item(198, 4).Value = 3.015625

$ws.Cells.Item(199, 1).Value = 3.6
$ws.Cells.Item(199, 2).Value = 0.04000000000000001
$ws.Cells.Item(199, 3).Value = 7.816614866256714
$ws.Cells.Item(199, 4).Value = 3.609375

$ws.Cells.Item(200, 1).Value = 3.6
$ws.Cells.Item(200, 2).Value = 0.04000000000000001
$ws.Cells.Item(200, 3).Value = 4.427245378494263
$ws.Cells.Item(200, 4).Value = 2.328125

$ws.Cells.Item(201, 1).Value = 59.83
$ws.Cells.Item(201, 2).Value = 0.04000000000000001
$ws.Cells.Item(201, 3).Value = 8.842975854873657
$ws.Cells.Item(201, 4).Value = 3.484375

$ws.Cells.Item(202, 1).Value = 39.06
$ws.Cells.Item(202, 2).Value = 0.04000000000000001
$ws.Cells.Item(202, 3).Value = 103.8992028236389
$ws.Cells.Item(202, 4).Value = 4.765625

$ws.Cells.Item(203, 1).Value = 39.06
$ws.Cells.Item(203, 2).Value = 0.04000000000000001
$ws.Cells.Item(203, 3).Value = 7.417251348495483
$ws.Cells.Item(203, 4).Value = 2.984375

$ws.Cells.Item(204, 1).Value = 39.78
$ws.Cells.Item(204, 2).Value = 0.16
$ws.Cells.Item(204, 3).Value = 22.68746423721313
$ws.Cells.Item(204, 4).Value = 13.84375

$ws.Cells.Item(205, 1).Value = 39.78
$ws.Cells.Item(205, 2).Value = 0.16
$ws.Cells.Item(205, 3).Value = 43.00482630729675
$ws.Cells.Item(205, 4).Value = 18.4375

$ws.Cells.Item(206, 1).Value = 39.06
$ws.Cells.Item(206, 2).Value = 0.04000000000000001
$ws.Cells.Item(206, 3).Value = 6.663230180740356
$ws.Cells.Item(206, 4).Value = 3.4375

$ws.Cells.Item(207, 1).Value = 39.06
$ws.Cells.Item(207, 2).Value = 0.04000000000000001
$ws.Cells.Item(207, 3).Value = 7.120179891586304
$ws.Cells.Item(207, 4).Value = 3.390625

$ws.Cells.Item(208, 1).Value = 39.78
$ws.Cells.Item(208, 2).Value = 0.16
$ws.Cells.Item(208, 3).Value = 22.58138370513916
$ws.Cells.Item(208, 4).Value = 14.84375

# Drop the now-stale row 209 (data shifted out of the cleaned table).
$ws.Rows.Item(209).Delete()
